# Cambio en llave primaria de NIT por codigoEntidad y se permite la carga
# multiples veces de archivo entidades.
#
# A new row is inserted right above the current row 13 (NIT 900226715 /
# codigo EPS042). The new row duplicates that entity's data but with a new
# "codigo" value (ESSC24), showing the same NIT can now appear more than
# once with a different codigoEntidad.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codigo" (B) and "razonEntidad" (C) columns now get their own
# explicit column-width entry (still the same visual width as the sheet
# default, just marked as a custom width instead of inheriting it).
$ws.Columns.Item(2).ColumnWidth = 21
$ws.Columns.Item(3).ColumnWidth = 21

# Insert a new blank row at row 13, pushing the existing row 13 (and
# everything below it) down by one.
$ws.Rows.Item(13).Insert()

# Fill the new row 13 with a copy of the entity that used to sit there
# (NIT 900226715), but with the new codigoEntidad "ESSC24".
$ws.Cells.Item(13, 1).Value = 900226715
$ws.Cells.Item(13, 2).Value = "ESSC24"
$ws.Cells.Item(13, 3).Value = "contribuciones"
$ws.Cells.Item(13, 4).Value = "COOSALUD "
$ws.Cells.Item(13, 5).Value = "SALUD"
$ws.Cells.Item(13, 6).Value = "A010102002"
$ws.Cells.Item(13, 7).Value = 73
$ws.Cells.Item(13, 8).Value = 23001010102

# Move the active selection to match the saved state of the workbook.
$ws.Range("C26").Select()
